$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 335; this shifts the existing rows
# 335-353 down to 337-355 and carries the row-335 formatting (incl. the
# date style on column D) onto the freshly inserted rows.
$ws.Rows("335:336").Insert()

# Row 335: new "Doctor Davis" / "Especial" record
$ws.Cells.Item(335, 1).Value = 5
$ws.Cells.Item(335, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(335, 3).Value = "Maule"
$ws.Cells.Item(335, 4).Value = 44610
$ws.Cells.Item(335, 5).Value = 7
$ws.Cells.Item(335, 6).Value = "Fruta"
$ws.Cells.Item(335, 7).Value = 100103
$ws.Cells.Item(335, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(335, 9).Value = 100103004
$ws.Cells.Item(335, 10).Value = "Durazno"
$ws.Cells.Item(335, 11).Value = "Doctor Davis"
$ws.Cells.Item(335, 12).Value = "Especial"
$ws.Cells.Item(335, 13).Value = 300
$ws.Cells.Item(335, 14).Value = 13000
$ws.Cells.Item(335, 15).Value = 13000
$ws.Cells.Item(335, 16).Value = 13000
$ws.Cells.Item(335, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(335, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(335, 19).Value = 867
$ws.Cells.Item(335, 20).Value = 15

# Row 336: new "Doctor Davis" / "Extra (doble especial)" record
$ws.Cells.Item(336, 1).Value = 5
$ws.Cells.Item(336, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(336, 3).Value = "Maule"
$ws.Cells.Item(336, 4).Value = 44610
$ws.Cells.Item(336, 5).Value = 7
$ws.Cells.Item(336, 6).Value = "Fruta"
$ws.Cells.Item(336, 7).Value = 100103
$ws.Cells.Item(336, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(336, 9).Value = 100103004
$ws.Cells.Item(336, 10).Value = "Durazno"
$ws.Cells.Item(336, 11).Value = "Doctor Davis"
$ws.Cells.Item(336, 12).Value = "Extra (doble especial)"
$ws.Cells.Item(336, 13).Value = 400
$ws.Cells.Item(336, 14).Value = 15000
$ws.Cells.Item(336, 15).Value = 15000
$ws.Cells.Item(336, 16).Value = 15000
$ws.Cells.Item(336, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(336, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(336, 19).Value = 1000
$ws.Cells.Item(336, 20).Value = 15
